$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the hyperlink that currently lives on A52 (the Eurostat database URL).
$ws.Range("A52").Hyperlinks.Delete()

# 2. The "SBS Main Indicators..." label moves down one row (A51 -> A52) and the
#    hyperlinked URL text moves further down to make room for a new blank
#    separator row. A51 becomes blank.
$ws.Range("A51").Value = ""
$ws.Range("A51").Style = "source"

$ws.Range("A52").Value = "SBS Main Indicators, Annual enterprise statistics by size class for special aggregates of activities (NACE Rev. 2)"
$ws.Range("A52").Style = "source"

# 3. Insert a new blank row at row 53 (pushes old row53.. down by one).
$ws.Rows.Item(53).Insert()

$ws.Range("A53").Value = ""
$ws.Range("A53").Style = "source"

$ws.Range("A54").Value = "http://epp.eurostat.ec.europa.eu/portal/page/portal/european_business/data/database"
$ws.Range("A54").Style = "source"

# 4. The closing citation line (now at A58 after the insert) is replaced with
#    the short "SBS Eurostat" label instead of the old long citation text.
$ws.Range("A58").Value = "SBS Eurostat"
$ws.Range("A58").Style = "source"
